$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 527.5714
$ws.Range("I38").Value = 32.166668
$ws.Range("J38").Value = 3500
$ws.Range("K38").Value = 96.500004
$ws.Range("L38").Value = 10500
$ws.Range("M38").Value = 275.499996
$ws.Range("N38").Value = -11244

$ws.Range("H116").Value = 5992.2
$ws.Range("I116").Value = 5989.5
$ws.Range("J116").Value = 5994
$ws.Range("K116").Value = 5989.5
$ws.Range("L116").Value = 5994
$ws.Range("M116").Value = -2547.5
$ws.Range("N116").Value = -12878

$ws.Range("H125").Value = 3803.4119
$ws.Range("I125").Value = 2635.3076
$ws.Range("K125").Value = 23717.7684
$ws.Range("M125").Value = -21257.7684

$ws.Range("H132").Value = 2304.4119
$ws.Range("I132").Value = 2304.4119
$ws.Range("K132").Value = 6913.2357
$ws.Range("M132").Value = -4383.2357

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1564
$ws.Range("I2").Value = 943.44446
$ws.Range("J2").Value = 4356.5
$ws.Range("K2").Value = 943.44446
$ws.Range("L2").Value = 4356.5
$ws.Range("M2").Value = -830.44446
$ws.Range("N2").Value = -4582.5

$ws.Range("H61").Value = 3057.4285
$ws.Range("I61").Value = 2880.6
$ws.Range("J61").Value = 3499.5
$ws.Range("K61").Value = 2880.6
$ws.Range("L61").Value = 3499.5
$ws.Range("M61").Value = -2668.6
$ws.Range("N61").Value = -3923.5

$ws.Range("H74").Value = 1783.3572
$ws.Range("I74").Value = 1875.1538
$ws.Range("J74").Value = 590
$ws.Range("K74").Value = 1875.1538
$ws.Range("L74").Value = 590
$ws.Range("M74").Value = -1001.1538
$ws.Range("N74").Value = -2338

$ws.Range("H77").Value = 1783.3572
$ws.Range("I77").Value = 1875.1538
$ws.Range("J77").Value = 590
$ws.Range("K77").Value = 9375.769
$ws.Range("L77").Value = 2950
$ws.Range("M77").Value = -5007.769
$ws.Range("N77").Value = -11686

$ws.Range("H116").Value = 1564
$ws.Range("I116").Value = 943.44446
$ws.Range("J116").Value = 4356.5
$ws.Range("K116").Value = 943.44446
$ws.Range("L116").Value = 4356.5
$ws.Range("M116").Value = 1350.55554
$ws.Range("N116").Value = -8944.5

$ws.Range("H135").Value = 41000
$ws.Range("J135").Value = 41000
$ws.Range("L135").Value = 41000
$ws.Range("N135").Value = -51140

$ws.Range("H136").Value = 3057.4285
$ws.Range("I136").Value = 2880.6
$ws.Range("J136").Value = 3499.5
$ws.Range("K136").Value = 8641.799999999999
$ws.Range("L136").Value = 10498.5
$ws.Range("M136").Value = -6091.799999999999
$ws.Range("N136").Value = -15598.5

$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1564
$ws.Range("I3").Value = 943.44446
$ws.Range("J3").Value = 4356.5
$ws.Range("K3").Value = 943.44446
$ws.Range("L3").Value = 4356.5
$ws.Range("M3").Value = -829.44446
$ws.Range("N3").Value = -4584.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 482.875
$ws.Range("I19").Value = 155.38461
$ws.Range("J19").Value = 1902
$ws.Range("K19").Value = 155.38461
$ws.Range("L19").Value = 1902
$ws.Range("M19").Value = 14.61538999999999
$ws.Range("N19").Value = -2242

$ws.Range("H24").Value = 482.875
$ws.Range("I24").Value = 155.38461
$ws.Range("J24").Value = 1902
$ws.Range("K24").Value = 155.38461
$ws.Range("L24").Value = 1902
$ws.Range("M24").Value = 14.61538999999999
$ws.Range("N24").Value = -2242

$ws.Range("H31").Value = 2818.4783
$ws.Range("I31").Value = 2048.25
$ws.Range("J31").Value = 4579
$ws.Range("K31").Value = 2048.25
$ws.Range("L31").Value = 4579
$ws.Range("M31").Value = -1753.25
$ws.Range("N31").Value = -5169

$ws.Range("H34").Value = 2818.4783
$ws.Range("I34").Value = 2048.25
$ws.Range("J34").Value = 4579
$ws.Range("K34").Value = 2048.25
$ws.Range("L34").Value = 4579
$ws.Range("M34").Value = -1846.25
$ws.Range("N34").Value = -4983

$ws.Range("H58").Value = 2103.9092
$ws.Range("I58").Value = 2103.9092
$ws.Range("K58").Value = 2103.9092
$ws.Range("M58").Value = -1900.9092

$ws.Range("H122").Value = 760.375
$ws.Range("I122").Value = 748.8333
$ws.Range("J122").Value = 795
$ws.Range("K122").Value = 2246.4999
$ws.Range("L122").Value = 2385
$ws.Range("M122").Value = 203.5001000000002
$ws.Range("N122").Value = -7285

$ws.Range("H132").Value = 2077.8
$ws.Range("I132").Value = 1270.1428
$ws.Range("J132").Value = 3962.3333
$ws.Range("K132").Value = 3810.4284
$ws.Range("L132").Value = 11886.9999
$ws.Range("M132").Value = -1280.4284
$ws.Range("N132").Value = -16946.9999

$ws.Range("H136").Value = 2103.9092
$ws.Range("I136").Value = 2103.9092
$ws.Range("K136").Value = 6311.7276
$ws.Range("M136").Value = -3761.7276

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 93996
$ws.Range("J37").Value = 93996
$ws.Range("L37").Value = 281988
$ws.Range("N37").Value = -282212

$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()

$ws.Range("H122").Value = 705
$ws.Range("I122").Value = 766.3333
$ws.Range("J122").Value = 671.5454999999999
$ws.Range("K122").Value = 6896.9997
$ws.Range("L122").Value = 6043.9095
$ws.Range("M122").Value = -4446.9997
$ws.Range("N122").Value = -10943.9095

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2305.5
$ws.Range("I132").Value = 2305.5
$ws.Range("K132").Value = 6916.5
$ws.Range("M132").Value = -4386.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 2599.4
$ws.Range("I4").Value = 2599.4
$ws.Range("K4").Value = 2599.4
$ws.Range("M4").Value = -2486.4

$ws.Range("H7").Value = 3978.1428
$ws.Range("I7").Value = 3978.1428
$ws.Range("K7").Value = 3978.1428
$ws.Range("M7").Value = -3866.1428

$ws.Range("H22").Value = 2064.1428
$ws.Range("I22").Value = 2033.1666
$ws.Range("J22").Value = 2250
$ws.Range("K22").Value = 2033.1666
$ws.Range("L22").Value = 2250
$ws.Range("M22").Value = -1738.1666
$ws.Range("N22").Value = -2840

$ws.Range("H27").Value = 2064.1428
$ws.Range("I27").Value = 2033.1666
$ws.Range("J27").Value = 2250
$ws.Range("K27").Value = 2033.1666
$ws.Range("L27").Value = 2250
$ws.Range("M27").Value = -1926.1666
$ws.Range("N27").Value = -2464

$ws.Range("H28").Value = 2599.4
$ws.Range("I28").Value = 2599.4
$ws.Range("K28").Value = 2599.4
$ws.Range("M28").Value = -2367.4

$ws.Range("H37").Value = 2599.4
$ws.Range("I37").Value = 2599.4
$ws.Range("K37").Value = 2599.4
$ws.Range("M37").Value = -2492.4

$ws.Range("H61").Value = 701.3333
$ws.Range("I61").Value = 197.66667
$ws.Range("K61").Value = 197.66667
$ws.Range("M61").Value = 4.333329999999989

$ws.Range("H113").Value = 701.3333
$ws.Range("I113").Value = 197.66667
$ws.Range("K113").Value = 197.66667
$ws.Range("M113").Value = 1972.33333

$ws.Range("H126").Value = 3978.1428
$ws.Range("I126").Value = 3978.1428
$ws.Range("K126").Value = 11934.4284
$ws.Range("M126").Value = -9464.428400000001

$ws.Range("H132").Value = 4614.2856
$ws.Range("I132").Value = 2512.5
$ws.Range("K132").Value = 7537.5
$ws.Range("M132").Value = -5007.5

$ws.Range("H136").Value = 2341.4167
$ws.Range("I136").Value = 2209.7
$ws.Range("K136").Value = 6629.099999999999
$ws.Range("M136").Value = -4079.099999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 20750.75
$ws.Range("I4").Value = 60000
$ws.Range("J4").Value = 7667.6665
$ws.Range("K4").Value = 60000
$ws.Range("L4").Value = 7667.6665
$ws.Range("M4").Value = -59887
$ws.Range("N4").Value = -7893.6665

$ws.Range("H23").Value = 15371.429
$ws.Range("I23").Value = 15371.429
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 15371.429
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -15142.429
$ws.Range("N23").ClearContents()

$ws.Range("H46").Value = 43333.332
$ws.Range("J46").Value = 43333.332
$ws.Range("L46").Value = 43333.332
$ws.Range("N46").Value = -43795.332

$ws.Range("H107").Value = 362.16666
$ws.Range("I107").Value = 293.25
$ws.Range("K107").Value = 879.75
$ws.Range("M107").Value = 1040.25

$ws.Range("H132").Value = 1530.8572
$ws.Range("I132").Value = 1365.0625
$ws.Range("J132").Value = 2061.4
$ws.Range("K132").Value = 4095.1875
$ws.Range("L132").Value = 6184.200000000001
$ws.Range("M132").Value = -1565.1875
$ws.Range("N132").Value = -11244.2

$ws.Range("H134").Value = 43333.332
$ws.Range("J134").Value = 43333.332
$ws.Range("L134").Value = 129999.996
$ws.Range("N134").Value = -135069.996
